# Regenerate the "K" column (formerly Strike#) values in column G.
# New K values were recalculated (std/mean, s_vals) for rows 2-46.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 2
    8  = 2
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 2
    16 = 3
    17 = 1
    18 = 1
    19 = 0
    20 = 2
    21 = 0
    22 = 0
    23 = 1
    24 = 1
    25 = 3
    26 = 0
    27 = 3
    28 = 3
    29 = 1
    30 = 1
    31 = 1
    32 = 2
    33 = 2
    34 = 1
    35 = 1
    36 = 1
    37 = 1
    38 = 2
    39 = 3
    40 = 3
    41 = 2
    42 = 4
    43 = 6
    44 = 5
    45 = 6
    46 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
